$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$s5 = $wb.Worksheets.Add($null, $lastSheet)
$s5.Name = "Cool AFDW per day Boot"

$s6 = $wb.Worksheets.Add($null, $s5)
$s6.Name = "Warm AFDW per day Boot"

$s7 = $wb.Worksheets.Add($null, $s6)
$s7.Name = "cool AFDW per day glm"

$s8 = $wb.Worksheets.Add($null, $s7)
$s8.Name = "warm  AFDW per day glm"

# --- Cool AFDW per day Boot ---
$s5.Columns.Item(1).ColumnWidth = 14.76
$s5.Columns.Item(2).ColumnWidth = 17.76
$s5.Columns.Item(3).ColumnWidth = 10.76
$s5.Columns.Item(4).ColumnWidth = 12.76
$s5.Columns.Item(5).ColumnWidth = 10.76
$s5.Columns.Item(6).ColumnWidth = 10.76
$s5.Range("A1").Value = "Treatment"
$s5.Range("B1").Value = "N"
$s5.Range("C1").Value = "mean"
$s5.Range("D1").Value = "boot_estimate"
$s5.Range("E1").Value = "lowerCI"
$s5.Range("F1").Value = "upperCI"
$s5.Range("A2").Value = "Open pizzas"
$s5.Range("B2").Value = "Ambient nutrients"
$s5.Range("C2").Value = 0.932833333333333
$s5.Range("D2").Value = 0.932833333333333
$s5.Range("E2").Value = 0.00248414115646259
$s5.Range("F2").Value = 0.0119289186507937
$s5.Range("A3").Value = "Open pizzas"
$s5.Range("B3").Value = "Nutrient additions"
$s5.Range("C3").Value = 3.08633333333333
$s5.Range("D3").Value = 3.08633333333333
$s5.Range("E3").Value = 0.00578662131519274
$s5.Range("F3").Value = 0.054271587326907
$s5.Range("A4").Value = "Fish exclusions"
$s5.Range("B4").Value = "Ambient nutrients"
$s5.Range("C4").Value = 10.7585714285714
$s5.Range("D4").Value = 10.7585714285714
$s5.Range("E4").Value = 0.0313214285714312
$s5.Range("F4").Value = 0.147286458333333
$s5.Range("A5").Value = "Fish exclusions"
$s5.Range("B5").Value = "Nutrient additions"
$s5.Range("C5").Value = 6.212875
$s5.Range("D5").Value = 6.212875
$s5.Range("E5").Value = 0.0318099702380951
$s5.Range("F5").Value = 0.0651238945578228
$s5.Range("A6").Value = "Full exclusions"
$s5.Range("B6").Value = "Ambient nutrients"
$s5.Range("C6").Value = 8.85225
$s5.Range("D6").Value = 8.85225
$s5.Range("E6").Value = 0.0236548752834467
$s5.Range("F6").Value = 0.120829828042328
$s5.Range("A7").Value = "Full exclusions"
$s5.Range("B7").Value = "Nutrient additions"
$s5.Range("C7").Value = 5.329125
$s5.Range("D7").Value = 5.329125
$s5.Range("E7").Value = 0.0256150793650794
$s5.Range("F7").Value = 0.0680006944444445

# --- Warm AFDW per day Boot ---
$s6.Columns.Item(1).ColumnWidth = 14.76
$s6.Columns.Item(2).ColumnWidth = 17.76
$s6.Columns.Item(3).ColumnWidth = 10.76
$s6.Columns.Item(4).ColumnWidth = 12.76
$s6.Columns.Item(5).ColumnWidth = 10.76
$s6.Columns.Item(6).ColumnWidth = 10.76
$s6.Range("A1").Value = "Treatment"
$s6.Range("B1").Value = "N"
$s6.Range("C1").Value = "mean"
$s6.Range("D1").Value = "boot_estimate"
$s6.Range("E1").Value = "lowerCI"
$s6.Range("F1").Value = "upperCI"
$s6.Range("A2").Value = "Open pizzas"
$s6.Range("B2").Value = "Ambient nutrients"
$s6.Range("C2").Value = 7.62014285714286
$s6.Range("D2").Value = 7.62014285714286
$s6.Range("E2").Value = 0.0184600549471694
$s6.Range("F2").Value = 0.378997692439591
$s6.Range("A3").Value = "Open pizzas"
$s6.Range("B3").Value = "Nutrient additions"
$s6.Range("C3").Value = 1.49325
$s6.Range("D3").Value = 1.49325
$s6.Range("E3").Value = 0.0156740551845135
$s6.Range("F3").Value = 0.0511373188405798
$s6.Range("A4").Value = "Full exclusions"
$s6.Range("B4").Value = "Ambient nutrients"
$s6.Range("C4").Value = 2.853125
$s6.Range("D4").Value = 2.853125
$s6.Range("E4").Value = 0.0286426921583851
$s6.Range("F4").Value = 0.114610688405797
$s6.Range("A5").Value = "Full exclusions"
$s6.Range("B5").Value = "Nutrient additions"
$s6.Range("C5").Value = 29.5645714285714
$s6.Range("D5").Value = 29.5645714285714
$s6.Range("E5").Value = 0.0672725414078676
$s6.Range("F5").Value = 1.41082527173913

# --- cool AFDW per day glm ---
$s7.Columns.Item(1).ColumnWidth = 43.76
$s7.Columns.Item(2).ColumnWidth = 10.76
$s7.Columns.Item(3).ColumnWidth = 10.76
$s7.Columns.Item(4).ColumnWidth = 10.76
$s7.Range("A1").Value = "Coefficient"
$s7.Range("B1").Value = "Estimate"
$s7.Range("C1").Value = "Std. Error"
$s7.Range("D1").Value = "t value"
$s7.Range("E1").Value = "P value"
$s7.Range("A2").Value = "(Intercept)"
$s7.Range("B2").Value = -4.90581063627031
$s7.Range("C2").Value = 0.417417161113429
$s7.Range("D2").Value = -11.7527765825066
$s7.Range("E2").Value = 0.0000000000000150328762006684
$s7.Range("A3").Value = "TreatmentFish exclusions"
$s7.Range("B3").Value = 2.44523150837728
$s7.Range("C3").Value = 0.568843704757533
$s7.Range("D3").Value = 4.29859992811126
$s7.Range("E3").Value = 0.000107048816386244
$s7.Range("A4").Value = "TreatmentFull exclusions"
$s7.Range("B4").Value = 2.25020039331424
$s7.Range("C4").Value = 0.552191000638357
$s7.Range("D4").Value = 4.07503996029075
$s7.Range("E4").Value = 0.000211970319649157
$s7.Range("A5").Value = "NNutrient additions"
$s7.Range("B5").Value = 1.19651250974114
$s7.Range("C5").Value = 0.538883237804496
$s7.Range("D5").Value = 2.22035577617136
$s7.Range("E5").Value = 0.0321243217145331
$s7.Range("A6").Value = "TreatmentFish exclusions:NNutrient additions"
$s7.Range("B6").Value = -1.74559153702552
$s7.Range("C6").Value = 0.755260903814658
$s7.Range("D6").Value = -2.31124307932387
$s7.Range("E6").Value = 0.0260543194564904
$s7.Range("A7").Value = "TreatmentFull exclusions:NNutrient additions"
$s7.Range("B7").Value = -1.70399711413855
$s7.Range("C7").Value = 0.742799282158141
$s7.Range("D7").Value = -2.29402094895371
$s7.Range("E7").Value = 0.0271187899854395

# --- warm  AFDW per day glm ---
$s8.Columns.Item(1).ColumnWidth = 43.76
$s8.Columns.Item(2).ColumnWidth = 10.76
$s8.Columns.Item(3).ColumnWidth = 10.76
$s8.Columns.Item(4).ColumnWidth = 10.76
$s8.Range("A1").Value = "Coefficient"
$s8.Range("B1").Value = "Estimate"
$s8.Range("C1").Value = "Std. Error"
$s8.Range("D1").Value = "t value"
$s8.Range("E1").Value = "P value"
$s8.Range("A2").Value = "(Intercept)"
$s8.Range("B2").Value = -1.79784627930817
$s8.Range("C2").Value = 0.477744116713901
$s8.Range("D2").Value = -3.76319920310986
$s8.Range("E2").Value = 0.000864300912356545
$s8.Range("A3").Value = "TreatmentFull exclusions"
$s8.Range("B3").Value = -0.982380232379726
$s8.Range("C3").Value = 0.65417807359896
$s8.Range("D3").Value = -1.50170155807143
$s8.Range("E3").Value = 0.145222728691909
$s8.Range("A4").Value = "NNutrient additions"
$s8.Range("B4").Value = -1.62984016455064
$s8.Range("C4").Value = 0.65417807359896
$s8.Range("D4").Value = -2.49143196681001
$s8.Range("E4").Value = 0.0194303213505482
$s8.Range("A5").Value = "TreatmentFull exclusions:NNutrient additions"
$s8.Range("B5").Value = 3.96800201301227
$s8.Range("C5").Value = 0.925147503890755
$s8.Range("D5").Value = 4.28904795864944
$s8.Range("E5").Value = 0.000219269141785231

# Restore the original active sheet/tab (new sheets were appended at
# the end and Add() shifts selection onto them as it goes).
$wb.Worksheets.Item(1).Activate()

